$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows in column C ("Started") whose value toggles between "Yes" and "No"
$ws.Range("C4").Value = "No"
$ws.Range("C5").Value = "No"
$ws.Range("C7").Value = "Yes"
$ws.Range("C8").Value = "Yes"
$ws.Range("C17").Value = "No"
$ws.Range("C23").Value = "Yes"
$ws.Range("C28").Value = "No"
$ws.Range("C37").Value = "Yes"
$ws.Range("C42").Value = "Yes"
$ws.Range("C43").Value = "No"
$ws.Range("C50").Value = "No"
$ws.Range("C53").Value = "Yes"
$ws.Range("C56").Value = "No"
$ws.Range("C60").Value = "Yes"
$ws.Range("C76").Value = "Yes"
$ws.Range("C79").Value = "No"

# Reflect the final selected cell as recorded in the saved workbook
$ws.Range("C38").Select()
